$d = $word.ActiveDocument

# Paragraph 1: select the whole paragraph's text range and apply strikethrough
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.StrikeThrough = 1

# Insert "[Done] " at the very start of the document (before existing text), without strikethrough
$start = $d.Range(0, 0)
$start.InsertBefore("[Done] ")
$start.Font.StrikeThrough = 0

# Update the OLEObject's ObjectID
$d.Content.Find.Execute("_1567950625", $false, $false, $false, $false, $false, $true, 1, $false, "_1568190480", 2)
